$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores runs/balls/fours/sixes as text-typed cells (numbers
# entered/stored as strings). Force Text number format first so the
# values we write stay text (t="str"/shared-string) instead of being
# reinterpreted as numeric cells.
$ws.Range("C3:E4").NumberFormat = "@"

# Row 3 (Sarfaraz Khan, Kings XI Punjab): runs 7->12, balls 8->12, fours 1->2
$ws.Range("C3").Value = "12"
$ws.Range("D3").Value = "12"
$ws.Range("E3").Value = "2"

# Row 4 (Sarfaraz Khan, Kings XI Punjab): runs 12->7, balls 12->8, fours 2->1
$ws.Range("C4").Value = "7"
$ws.Range("D4").Value = "8"
$ws.Range("E4").Value = "1"
